# Auto-generated edit script: update cached Leve profit values across 8 Sephirot_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 3245.5833
$ws.Cells.Item(17, 10).Value = 3245.5833
$ws.Cells.Item(17, 12).Value = 9736.749899999999
$ws.Cells.Item(17, 14).Value = -10072.7499

$ws.Cells.Item(19, 8).Value = 808
$ws.Cells.Item(19, 9).Value = 616.6667
$ws.Cells.Item(19, 11).Value = 616.6667
$ws.Cells.Item(19, 13).Value = -441.6667

$ws.Cells.Item(98, 8).Value = 398.5
$ws.Cells.Item(98, 9).Value = 398.5
$ws.Cells.Item(98, 11).Value = 398.5
$ws.Cells.Item(98, 13).Value = 1099.5

$ws.Cells.Item(122, 8).Value = 398.5
$ws.Cells.Item(122, 9).Value = 398.5
$ws.Cells.Item(122, 11).Value = 1195.5
$ws.Cells.Item(122, 13).Value = 1254.5

$ws.Cells.Item(132, 8).Value = 2979.8
$ws.Cells.Item(132, 9).Value = 2974.75
$ws.Cells.Item(132, 10).Value = 3000
$ws.Cells.Item(132, 11).Value = 8924.25
$ws.Cells.Item(132, 12).Value = 9000
$ws.Cells.Item(132, 13).Value = -6394.25
$ws.Cells.Item(132, 14).Value = -14060

$ws.Cells.Item(137, 8).Value = 2825.5715
$ws.Cells.Item(137, 10).Value = 3499
$ws.Cells.Item(137, 12).Value = 10497
$ws.Cells.Item(137, 14).Value = -15597

$ws.Cells.Item(138, 8).Value = 5475.387
$ws.Cells.Item(138, 10).Value = 5525.448
$ws.Cells.Item(138, 12).Value = 16576.344
$ws.Cells.Item(138, 14).Value = -26856.344

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 499.5
$ws.Cells.Item(2, 9).Value = 499.5
$ws.Cells.Item(2, 11).Value = 499.5
$ws.Cells.Item(2, 13).Value = -386.5

$ws.Cells.Item(32, 8).Value = 4615.0347
$ws.Cells.Item(32, 9).Value = 4615.0347
$ws.Cells.Item(32, 11).Value = 4615.0347
$ws.Cells.Item(32, 13).Value = -4328.0347

$ws.Cells.Item(61, 8).Value = 2689.3333
$ws.Cells.Item(61, 9).Value = 1309.909
$ws.Cells.Item(61, 10).Value = 4857
$ws.Cells.Item(61, 11).Value = 1309.909
$ws.Cells.Item(61, 12).Value = 4857
$ws.Cells.Item(61, 13).Value = -1097.909
$ws.Cells.Item(61, 14).Value = -5281

$ws.Cells.Item(116, 8).Value = 499.5
$ws.Cells.Item(116, 9).Value = 499.5
$ws.Cells.Item(116, 11).Value = 499.5
$ws.Cells.Item(116, 13).Value = 1794.5

$ws.Cells.Item(136, 8).Value = 2689.3333
$ws.Cells.Item(136, 9).Value = 1309.909
$ws.Cells.Item(136, 10).Value = 4857
$ws.Cells.Item(136, 11).Value = 3929.727
$ws.Cells.Item(136, 12).Value = 14571
$ws.Cells.Item(136, 13).Value = -1379.727
$ws.Cells.Item(136, 14).Value = -19671

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 499.5
$ws.Cells.Item(3, 9).Value = 499.5
$ws.Cells.Item(3, 11).Value = 499.5
$ws.Cells.Item(3, 13).Value = -385.5

$ws.Cells.Item(94, 8).Value = 1050
$ws.Cells.Item(94, 9).Value = 1050
$ws.Cells.Item(94, 10).Value = 0
$ws.Cells.Item(94, 11).Value = 1050
$ws.Cells.Item(94, 12).Value = 0
$ws.Cells.Item(94, 13).Value = -599
$ws.Cells.Item(94, 14).ClearContents()

$ws.Cells.Item(134, 8).Value = 15485.692
$ws.Cells.Item(134, 9).Value = 9210.362999999999
$ws.Cells.Item(134, 11).Value = 27631.089
$ws.Cells.Item(134, 13).Value = -25096.089

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(23, 8).Value = 49500
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 11).Value = 0
$ws.Cells.Item(23, 13).ClearContents()

$ws.Cells.Item(27, 8).Value = 49500
$ws.Cells.Item(27, 9).Value = 0
$ws.Cells.Item(27, 11).Value = 0
$ws.Cells.Item(27, 13).ClearContents()

$ws.Cells.Item(31, 8).Value = 3311.7058
$ws.Cells.Item(31, 9).Value = 1092.8572
$ws.Cells.Item(31, 10).Value = 4864.9
$ws.Cells.Item(31, 11).Value = 1092.8572
$ws.Cells.Item(31, 12).Value = 4864.9
$ws.Cells.Item(31, 13).Value = -797.8571999999999
$ws.Cells.Item(31, 14).Value = -5454.9

$ws.Cells.Item(34, 8).Value = 3311.7058
$ws.Cells.Item(34, 9).Value = 1092.8572
$ws.Cells.Item(34, 10).Value = 4864.9
$ws.Cells.Item(34, 11).Value = 1092.8572
$ws.Cells.Item(34, 12).Value = 4864.9
$ws.Cells.Item(34, 13).Value = -890.8571999999999
$ws.Cells.Item(34, 14).Value = -5268.9

$ws.Cells.Item(58, 8).Value = 1007.3333
$ws.Cells.Item(58, 9).Value = 1007.3333
$ws.Cells.Item(58, 10).Value = 0
$ws.Cells.Item(58, 11).Value = 1007.3333
$ws.Cells.Item(58, 12).Value = 0
$ws.Cells.Item(58, 13).Value = -804.3333
$ws.Cells.Item(58, 14).ClearContents()

$ws.Cells.Item(93, 8).Value = 0
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 2194.1177
$ws.Cells.Item(132, 9).Value = 1593.0714
$ws.Cells.Item(132, 11).Value = 4779.2142
$ws.Cells.Item(132, 13).Value = -2249.2142

$ws.Cells.Item(136, 8).Value = 1007.3333
$ws.Cells.Item(136, 9).Value = 1007.3333
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 3021.9999
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 13).Value = -471.9998999999998
$ws.Cells.Item(136, 14).ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 972.4286
$ws.Cells.Item(5, 9).Value = 967.8333
$ws.Cells.Item(5, 10).Value = 1000
$ws.Cells.Item(5, 11).Value = 2903.4999
$ws.Cells.Item(5, 12).Value = 3000
$ws.Cells.Item(5, 13).Value = -2791.4999
$ws.Cells.Item(5, 14).Value = -3224

$ws.Cells.Item(135, 8).Value = 972.4286
$ws.Cells.Item(135, 9).Value = 967.8333
$ws.Cells.Item(135, 10).Value = 1000
$ws.Cells.Item(135, 11).Value = 8710.4997
$ws.Cells.Item(135, 12).Value = 9000
$ws.Cells.Item(135, 13).Value = -6175.4997
$ws.Cells.Item(135, 14).Value = -14070

$ws.Cells.Item(138, 8).Value = 3188.3333
$ws.Cells.Item(138, 9).Value = 3026
$ws.Cells.Item(138, 11).Value = 9078
$ws.Cells.Item(138, 13).Value = -3938

$ws.Cells.Item(139, 8).Value = 1971
$ws.Cells.Item(139, 9).Value = 1971
$ws.Cells.Item(139, 11).Value = 5913
$ws.Cells.Item(139, 13).Value = -773

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 0
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 13).ClearContents()

$ws.Cells.Item(122, 8).Value = 8042.1
$ws.Cells.Item(122, 9).Value = 9301.625
$ws.Cells.Item(122, 10).Value = 3004
$ws.Cells.Item(122, 11).Value = 27904.875
$ws.Cells.Item(122, 12).Value = 9012
$ws.Cells.Item(122, 13).Value = -25454.875
$ws.Cells.Item(122, 14).Value = -13912

$ws.Cells.Item(132, 8).Value = 4166.3335
$ws.Cells.Item(132, 9).Value = 3750
$ws.Cells.Item(132, 11).Value = 11250
$ws.Cells.Item(132, 13).Value = -8720

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2000
$ws.Cells.Item(7, 9).Value = 2000
$ws.Cells.Item(7, 11).Value = 2000
$ws.Cells.Item(7, 13).Value = -1888

$ws.Cells.Item(126, 8).Value = 2000
$ws.Cells.Item(126, 9).Value = 2000
$ws.Cells.Item(126, 11).Value = 6000
$ws.Cells.Item(126, 13).Value = -3530

$ws.Cells.Item(136, 8).Value = 0
$ws.Cells.Item(136, 9).Value = 0
$ws.Cells.Item(136, 11).Value = 0
$ws.Cells.Item(136, 13).ClearContents()

$ws.Cells.Item(141, 8).Value = 62497.5
$ws.Cells.Item(141, 9).Value = 75000
$ws.Cells.Item(141, 10).Value = 49995
$ws.Cells.Item(141, 11).Value = 75000
$ws.Cells.Item(141, 12).Value = 49995
$ws.Cells.Item(141, 13).Value = -69820
$ws.Cells.Item(141, 14).Value = -60355

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 1615
$ws.Cells.Item(113, 9).Value = 1295
$ws.Cells.Item(113, 11).Value = 3885
$ws.Cells.Item(113, 13).Value = -1715

$ws.Cells.Item(132, 8).Value = 2747.238
$ws.Cells.Item(132, 9).Value = 2236.2666
$ws.Cells.Item(132, 10).Value = 4024.6667
$ws.Cells.Item(132, 11).Value = 6708.7998
$ws.Cells.Item(132, 12).Value = 12074.0001
$ws.Cells.Item(132, 13).Value = -4178.7998
$ws.Cells.Item(132, 14).Value = -17134.0001

$ws.Cells.Item(136, 8).Value = 816.35
$ws.Cells.Item(136, 9).Value = 789.5
$ws.Cells.Item(136, 11).Value = 2368.5
$ws.Cells.Item(136, 13).Value = 181.5
